# Applies the "456a3b4" gh-pages data refresh to 江西-漫展信息.xlsx
#
# Sheet "展览" (Exhibition) rows identified by row index on that sheet:
#   G2  : 258            -> "不可售"   (ticket became unsellable / no longer numeric)
#   F3  : 5615 -> 5620    (南昌·CM03动漫游戏博览会 - want-to-go count)
#   F6  : 930  -> 931     (南昌·第四届龙年动漫展)
#   F8  : 2513 -> 2516    (南昌·Sunflower Garden动漫游戏展)
#   F9  : 82   -> 83      (南昌·第一届哥布林动漫游戏展)
#   F10 : 141  -> 143     (赣州·卡尼动漫展)
#   F12 : 78   -> 79      (南昌·Aud中秋动漫嘉年华)
#   F13 : 14   -> 17      (九江·星梦次元XACD动漫游戏博览会国庆盛典)
#   G13 : 44.9 -> 55      (same row, min ticket price)
#   F14 : 2352 -> 2359    (南昌·萌卡动漫展)
#   F15 : 330  -> 340     (江西·JMG（广电）第二届UP动漫游戏博览会)
#
# Sheet "全部类型" (All types) mirrors the same events at different row
# numbers (it interleaves rows from the other category sheets):
#   G2  : 258  -> "不可售"
#   F3  : 5615 -> 5620
#   F8  : 930  -> 931
#   F10 : 2513 -> 2516
#   F12 : 141  -> 143
#   F16 : 14   -> 17
#   G16 : 44.9 -> 55
#   F17 : 2352 -> 2359
#   F18 : 330  -> 340

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F3").Value = 5620
$ws1.Range("F6").Value = 931
$ws1.Range("F8").Value = 2516
$ws1.Range("F9").Value = 83
$ws1.Range("F10").Value = 143
$ws1.Range("F12").Value = 79
$ws1.Range("F13").Value = 17
$ws1.Range("G13").Value = 55
$ws1.Range("F14").Value = 2359
$ws1.Range("F15").Value = 340

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F3").Value = 5620
$ws4.Range("F8").Value = 931
$ws4.Range("F10").Value = 2516
$ws4.Range("F12").Value = 143
$ws4.Range("F16").Value = 17
$ws4.Range("G16").Value = 55
$ws4.Range("F17").Value = 2359
$ws4.Range("F18").Value = 340
